$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2308.8572
$ws.Range("I62").Value = 2323.6667
$ws.Range("K62").Value = 2323.6667
$ws.Range("M62").Value = -1699.6667
$ws.Range("H65").Value = 2308.8572
$ws.Range("I65").Value = 2323.6667
$ws.Range("K65").Value = 11618.3335
$ws.Range("M65").Value = -8498.333500000001
$ws.Range("H68").Value = 38396
$ws.Range("J68").Value = 38396
$ws.Range("L68").Value = 38396
$ws.Range("N68").Value = -39894
$ws.Range("H71").Value = 38396
$ws.Range("J71").Value = 38396
$ws.Range("L71").Value = 115188
$ws.Range("N71").Value = -122676
$ws.Range("H129").Value = 860.0645
$ws.Range("I129").Value = 325
$ws.Range("J129").Value = 939.3333
$ws.Range("K129").Value = 975
$ws.Range("L129").Value = 2817.9999
$ws.Range("M129").Value = 4025
$ws.Range("N129").Value = -12817.9999
$ws.Range("H137").Value = 2351.9744
$ws.Range("I137").Value = 1897.24
$ws.Range("J137").Value = 3164
$ws.Range("K137").Value = 5691.72
$ws.Range("L137").Value = 9492
$ws.Range("M137").Value = -3141.72
$ws.Range("N137").Value = -14592
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1223.25
$ws.Range("I2").Value = 1441.2
$ws.Range("J2").Value = 1005.3
$ws.Range("K2").Value = 1441.2
$ws.Range("L2").Value = 1005.3
$ws.Range("M2").Value = -1328.2
$ws.Range("N2").Value = -1231.3
$ws.Range("H61").Value = 6218.3335
$ws.Range("I61").Value = 4324.4604
$ws.Range("J61").Value = 16161.167
$ws.Range("K61").Value = 4324.4604
$ws.Range("L61").Value = 16161.167
$ws.Range("M61").Value = -4112.4604
$ws.Range("N61").Value = -16585.167
$ws.Range("H74").Value = 4392.5137
$ws.Range("I74").Value = 1771.3715
$ws.Range("J74").Value = 50262.5
$ws.Range("K74").Value = 1771.3715
$ws.Range("L74").Value = 50262.5
$ws.Range("M74").Value = -897.3715
$ws.Range("N74").Value = -52010.5
$ws.Range("H77").Value = 4392.5137
$ws.Range("I77").Value = 1771.3715
$ws.Range("J77").Value = 50262.5
$ws.Range("K77").Value = 8856.8575
$ws.Range("L77").Value = 251312.5
$ws.Range("M77").Value = -4488.8575
$ws.Range("N77").Value = -260048.5
$ws.Range("H110").Value = 1283.6666
$ws.Range("I110").Value = 864.7143
$ws.Range("K110").Value = 864.7143
$ws.Range("M110").Value = 1180.2857
$ws.Range("H116").Value = 1223.25
$ws.Range("I116").Value = 1441.2
$ws.Range("J116").Value = 1005.3
$ws.Range("K116").Value = 1441.2
$ws.Range("L116").Value = 1005.3
$ws.Range("M116").Value = 852.8
$ws.Range("N116").Value = -5593.3
$ws.Range("H132").Value = 3168.4219
$ws.Range("I132").Value = 1020.04443
$ws.Range("J132").Value = 8256.684999999999
$ws.Range("K132").Value = 3060.13329
$ws.Range("L132").Value = 24770.055
$ws.Range("M132").Value = -530.1332900000002
$ws.Range("N132").Value = -29830.055
$ws.Range("H136").Value = 6218.3335
$ws.Range("I136").Value = 4324.4604
$ws.Range("J136").Value = 16161.167
$ws.Range("K136").Value = 12973.3812
$ws.Range("L136").Value = 48483.501
$ws.Range("M136").Value = -10423.3812
$ws.Range("N136").Value = -53583.501
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1223.25
$ws.Range("I3").Value = 1441.2
$ws.Range("J3").Value = 1005.3
$ws.Range("K3").Value = 1441.2
$ws.Range("L3").Value = 1005.3
$ws.Range("M3").Value = -1327.2
$ws.Range("N3").Value = -1233.3
$ws.Range("H99").Value = 1841.1
$ws.Range("I99").Value = 1685.7142
$ws.Range("K99").Value = 1685.7142
$ws.Range("M99").Value = -187.7141999999999
$ws.Range("H134").Value = 18030.656
$ws.Range("I134").Value = 1408.7347
$ws.Range("J134").Value = 85903.5
$ws.Range("K134").Value = 4226.2041
$ws.Range("L134").Value = 257710.5
$ws.Range("M134").Value = -1691.2041
$ws.Range("N134").Value = -262780.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2144.2263
$ws.Range("I31").Value = 1560.8837
$ws.Range("J31").Value = 4652.6
$ws.Range("K31").Value = 1560.8837
$ws.Range("L31").Value = 4652.6
$ws.Range("M31").Value = -1265.8837
$ws.Range("N31").Value = -5242.6
$ws.Range("H34").Value = 2144.2263
$ws.Range("I34").Value = 1560.8837
$ws.Range("J34").Value = 4652.6
$ws.Range("K34").Value = 1560.8837
$ws.Range("L34").Value = 4652.6
$ws.Range("M34").Value = -1358.8837
$ws.Range("N34").Value = -5056.6
$ws.Range("H94").Value = 1293.4286
$ws.Range("I94").Value = 1012
$ws.Range("K94").Value = 1012
$ws.Range("M94").Value = -561
$ws.Range("H129").Value = 45000
$ws.Range("J129").Value = 45000
$ws.Range("L129").Value = 45000
$ws.Range("N129").Value = -55000
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 862
$ws.Range("I122").Value = 562
$ws.Range("J122").Value = 1072
$ws.Range("K122").Value = 5058
$ws.Range("L122").Value = 9648
$ws.Range("M122").Value = -2608
$ws.Range("N122").Value = -14548
$ws.Range("H132").Value = 1700.1786
$ws.Range("I132").Value = 1841.4546
$ws.Range("J132").Value = 1608.7646
$ws.Range("K132").Value = 16573.0914
$ws.Range("L132").Value = 14478.8814
$ws.Range("M132").Value = -14043.0914
$ws.Range("N132").Value = -19538.8814
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 210001700
$ws.Range("I3").Value = 350000160
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 350000160
$ws.Range("L3").Value = 4000
$ws.Range("M3").Value = -350000044
$ws.Range("N3").Value = -4232
$ws.Range("H11").Value = 188336910
$ws.Range("I11").Value = 282500000
$ws.Range("J11").Value = 10750
$ws.Range("K11").Value = 282500000
$ws.Range("L11").Value = 10750
$ws.Range("M11").Value = -282499861
$ws.Range("N11").Value = -11028
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
$ws.Range("H24").Value = 422429.03
$ws.Range("I24").Value = 10007500
$ws.Range("J24").Value = 14553.681
$ws.Range("K24").Value = 10007500
$ws.Range("L24").Value = 14553.681
$ws.Range("M24").Value = -10007327
$ws.Range("N24").Value = -14899.681
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 3083.712
$ws.Range("I132").Value = 1430.6981
$ws.Range("J132").Value = 17685.334
$ws.Range("K132").Value = 4292.094300000001
$ws.Range("L132").Value = 53056.00199999999
$ws.Range("M132").Value = -1762.094300000001
$ws.Range("N132").Value = -58116.00199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1585.9474
$ws.Range("I16").Value = 1079.25
$ws.Range("K16").Value = 1079.25
$ws.Range("M16").Value = -909.25
$ws.Range("H20").Value = 3000
$ws.Range("J20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("N20").Value = -3452
$ws.Range("H46").Value = 1236
$ws.Range("J46").Value = 1295
$ws.Range("L46").Value = 1295
$ws.Range("N46").Value = -1671
$ws.Range("H100").Value = 5427.909
$ws.Range("I100").Value = 2925.75
$ws.Range("K100").Value = 2925.75
$ws.Range("M100").Value = -2384.75
$ws.Range("H136").Value = 2630.6086
$ws.Range("I136").Value = 1476.582
$ws.Range("J136").Value = 5723.4
$ws.Range("K136").Value = 4429.746
$ws.Range("L136").Value = 17170.2
$ws.Range("M136").Value = -1879.746
$ws.Range("N136").Value = -22270.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1356.4791
$ws.Range("I132").Value = 617.0606
$ws.Range("J132").Value = 2983.2
$ws.Range("K132").Value = 1851.1818
$ws.Range("L132").Value = 8949.599999999999
$ws.Range("M132").Value = 678.8181999999999
$ws.Range("N132").Value = -14009.6
